$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a running log table (A: date, B: weekday, C: hour, D: rank).
# Append the next day's entry as a new row right after the last used row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Column A holds the date as plain text (e.g. "2025/09/30"), not a real
# date value, so force text formatting before assigning it, then clear the
# formatting again so the cell ends up unstyled like its neighbours.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/09/30"
$ws.Cells.Item($newRow, 1).ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "火"
$ws.Cells.Item($newRow, 3).Value = 1
$ws.Cells.Item($newRow, 4).Value = 3
